# --------------------------------------------------------------------------
# Adds a "PO Forecast" sheet (Prophet-style forecast output) after the
# existing "Monthly Trend" sheet, and renames the "Requested quantity"
# headers on the two existing sheets to distinguish weekly vs. monthly
# quantity columns.
# --------------------------------------------------------------------------

$wb = $excel.ActiveWorkbook

# --- Rename the ambiguous "Requested quantity" header on each sheet ---
$wsWeekly  = $wb.Worksheets.Item("Weekly Quantity")
$wsWeekly.Range("B1").Value = "Weekly_PO_Qty"

$wsMonthly = $wb.Worksheets.Item("Monthly Trend")
$wsMonthly.Range("B1").Value = "Monthly_PO_Qty"

# --- Add the new "PO Forecast" sheet right after "Monthly Trend" ---
$wsForecast = $wb.Worksheets.Add($null, $wsMonthly)
$wsForecast.Name = "PO Forecast"

# Header row
$wsForecast.Cells.Item(1,1).Value = "ds"
$wsForecast.Cells.Item(1,2).Value = "PO_Forecast"
$wsForecast.Cells.Item(1,3).Value = "yhat_lower"
$wsForecast.Cells.Item(1,4).Value = "yhat_upper"

# Reuse the bold/bordered/centered header style from an existing sheet
$wsWeekly.Range("A1:B1").Copy()
$wsForecast.Range("A1:D1").PasteSpecial(-4122)

# Forecast data rows (ds, PO_Forecast, yhat_lower, yhat_upper)
$forecastData = New-Object 'object[,]' 57,4
$forecastData[0,0] = 44934.99999999999; $forecastData[0,1] = 1; $forecastData[0,2] = -3.442673459598333; $forecastData[0,3] = 6.599956146744357
$forecastData[1,0] = 44948.99999999999; $forecastData[1,1] = 2; $forecastData[1,2] = -3.686689277164881; $forecastData[1,3] = 6.372356363502315
$forecastData[2,0] = 44955.99999999999; $forecastData[2,1] = 2; $forecastData[2,2] = -3.58013610574377; $forecastData[2,3] = 6.873322667764203
$forecastData[3,0] = 44997.99999999999; $forecastData[3,1] = 2; $forecastData[3,2] = -3.303582101595801; $forecastData[3,3] = 6.983861502421952
$forecastData[4,0] = 45025.99999999999; $forecastData[4,1] = 2; $forecastData[4,2] = -2.785605845723687; $forecastData[4,3] = 7.41235104584276
$forecastData[5,0] = 45039.99999999999; $forecastData[5,1] = 2; $forecastData[5,2] = -3.416734014878271; $forecastData[5,3] = 7.341906270392292
$forecastData[6,0] = 45046.99999999999; $forecastData[6,1] = 2; $forecastData[6,2] = -2.512016257575424; $forecastData[6,3] = 7.598810651779669
$forecastData[7,0] = 45053.99999999999; $forecastData[7,1] = 2; $forecastData[7,2] = -2.470620411604102; $forecastData[7,3] = 7.82342180556673
$forecastData[8,0] = 45060.99999999999; $forecastData[8,1] = 3; $forecastData[8,2] = -2.565957673236944; $forecastData[8,3] = 7.520930386396699
$forecastData[9,0] = 45067.99999999999; $forecastData[9,1] = 3; $forecastData[9,2] = -2.451279403988729; $forecastData[9,3] = 7.963569701874232
$forecastData[10,0] = 45074.99999999999; $forecastData[10,1] = 3; $forecastData[10,2] = -2.009204542628888; $forecastData[10,3] = 7.88796285733345
$forecastData[11,0] = 45081.99999999999; $forecastData[11,1] = 3; $forecastData[11,2] = -2.36976616794369; $forecastData[11,3] = 7.752807282883344
$forecastData[12,0] = 45095.99999999999; $forecastData[12,1] = 3; $forecastData[12,2] = -2.313984932519999; $forecastData[12,3] = 8.018323317529033
$forecastData[13,0] = 45109.99999999999; $forecastData[13,1] = 3; $forecastData[13,2] = -1.908491841243086; $forecastData[13,3] = 8.034085037361947
$forecastData[14,0] = 45116.99999999999; $forecastData[14,1] = 3; $forecastData[14,2] = -2.320025762971188; $forecastData[14,3] = 8.220321635935113
$forecastData[15,0] = 45123.99999999999; $forecastData[15,1] = 3; $forecastData[15,2] = -2.017580395208201; $forecastData[15,3] = 7.973240125657296
$forecastData[16,0] = 45130.99999999999; $forecastData[16,1] = 3; $forecastData[16,2] = -2.177603795026601; $forecastData[16,3] = 7.899874928044508
$forecastData[17,0] = 45137.99999999999; $forecastData[17,1] = 3; $forecastData[17,2] = -1.78846142933331; $forecastData[17,3] = 8.170594464881683
$forecastData[18,0] = 45144.99999999999; $forecastData[18,1] = 3; $forecastData[18,2] = -1.786930345934127; $forecastData[18,3] = 8.392276117175113
$forecastData[19,0] = 45151.99999999999; $forecastData[19,1] = 3; $forecastData[19,2] = -1.59475170002141; $forecastData[19,3] = 8.431138202385036
$forecastData[20,0] = 45158.99999999999; $forecastData[20,1] = 3; $forecastData[20,2] = -1.63915411773082; $forecastData[20,3] = 8.98869462743138
$forecastData[21,0] = 45165.99999999999; $forecastData[21,1] = 3; $forecastData[21,2] = -1.665075313091744; $forecastData[21,3] = 8.567666083890986
$forecastData[22,0] = 45172.99999999999; $forecastData[22,1] = 3; $forecastData[22,2] = -1.622707881086851; $forecastData[22,3] = 8.535256768129312
$forecastData[23,0] = 45179.99999999999; $forecastData[23,1] = 4; $forecastData[23,2] = -1.782796223807383; $forecastData[23,3] = 8.650945248539733
$forecastData[24,0] = 45186.99999999999; $forecastData[24,1] = 4; $forecastData[24,2] = -1.779858814407107; $forecastData[24,3] = 8.19721164304179
$forecastData[25,0] = 45193.99999999999; $forecastData[25,1] = 4; $forecastData[25,2] = -1.511053253856098; $forecastData[25,3] = 8.723043854266795
$forecastData[26,0] = 45200.99999999999; $forecastData[26,1] = 4; $forecastData[26,2] = -1.178469688128716; $forecastData[26,3] = 8.652344007460146
$forecastData[27,0] = 45207.99999999999; $forecastData[27,1] = 4; $forecastData[27,2] = -1.238768937386436; $forecastData[27,3] = 8.686143677866097
$forecastData[28,0] = 45214.99999999999; $forecastData[28,1] = 4; $forecastData[28,2] = -1.500313257816676; $forecastData[28,3] = 8.727952421906433
$forecastData[29,0] = 45221.99999999999; $forecastData[29,1] = 4; $forecastData[29,2] = -1.164461234354023; $forecastData[29,3] = 9.047705206089661
$forecastData[30,0] = 45228.99999999999; $forecastData[30,1] = 4; $forecastData[30,2] = -1.231066931250805; $forecastData[30,3] = 9.170208418351082
$forecastData[31,0] = 45235.99999999999; $forecastData[31,1] = 4; $forecastData[31,2] = -1.087765484754121; $forecastData[31,3] = 9.035366845849198
$forecastData[32,0] = 45242.99999999999; $forecastData[32,1] = 4; $forecastData[32,2] = -0.6910849894070038; $forecastData[32,3] = 9.670108652132155
$forecastData[33,0] = 45249.99999999999; $forecastData[33,1] = 4; $forecastData[33,2] = -0.9307246107891231; $forecastData[33,3] = 9.205098701960097
$forecastData[34,0] = 45256.99999999999; $forecastData[34,1] = 4; $forecastData[34,2] = -0.7266792054768734; $forecastData[34,3] = 9.28928302508019
$forecastData[35,0] = 45263.99999999999; $forecastData[35,1] = 4; $forecastData[35,2] = -1.056107527854671; $forecastData[35,3] = 9.203924207587473
$forecastData[36,0] = 45270.99999999999; $forecastData[36,1] = 4; $forecastData[36,2] = -0.6094490109051158; $forecastData[36,3] = 9.280681424845165
$forecastData[37,0] = 45277.99999999999; $forecastData[37,1] = 4; $forecastData[37,2] = -0.7380868085263023; $forecastData[37,3] = 9.613437106624524
$forecastData[38,0] = 45298.99999999999; $forecastData[38,1] = 5; $forecastData[38,2] = -0.4630544303341775; $forecastData[38,3] = 9.64506730409528
$forecastData[39,0] = 45305.99999999999; $forecastData[39,1] = 5; $forecastData[39,2] = -0.3830757553183191; $forecastData[39,3] = 9.693707895125025
$forecastData[40,0] = 45312.99999999999; $forecastData[40,1] = 5; $forecastData[40,2] = -0.1409226718772111; $forecastData[40,3] = 9.649196446003723
$forecastData[41,0] = 45319.99999999999; $forecastData[41,1] = 5; $forecastData[41,2] = -0.1218522458064423; $forecastData[41,3] = 10.16707553724045
$forecastData[42,0] = 45333.99999999999; $forecastData[42,1] = 5; $forecastData[42,2] = -0.645640289575241; $forecastData[42,3] = 10.24583071021407
$forecastData[43,0] = 45340.99999999999; $forecastData[43,1] = 5; $forecastData[43,2] = -0.5798574231254536; $forecastData[43,3] = 10.21965881972277
$forecastData[44,0] = 45354.99999999999; $forecastData[44,1] = 5; $forecastData[44,2] = 0.1095868435505733; $forecastData[44,3] = 10.59138261382986
$forecastData[45,0] = 45361.99999999999; $forecastData[45,1] = 5; $forecastData[45,2] = -0.3232713438559027; $forecastData[45,3] = 10.2133217293707
$forecastData[46,0] = 45368.99999999999; $forecastData[46,1] = 5; $forecastData[46,2] = -0.1759572642207735; $forecastData[46,3] = 10.30703927287995
$forecastData[47,0] = 45375.99999999999; $forecastData[47,1] = 5; $forecastData[47,2] = -0.2615145163068194; $forecastData[47,3] = 10.465682984332
$forecastData[48,0] = 45382.99999999999; $forecastData[48,1] = 5; $forecastData[48,2] = 0.2224063447920794; $forecastData[48,3] = 10.37529500316323
$forecastData[49,0] = 45389.99999999999; $forecastData[49,1] = 5; $forecastData[49,2] = 0.2466513531635667; $forecastData[49,3] = 10.55085329556909
$forecastData[50,0] = 45396.99999999999; $forecastData[50,1] = 5; $forecastData[50,2] = -0.112715425529768; $forecastData[50,3] = 10.39399992682879
$forecastData[51,0] = 45403.99999999999; $forecastData[51,1] = 5; $forecastData[51,2] = 0.6245230644458247; $forecastData[51,3] = 10.58124242158025
$forecastData[52,0] = 45410.99999999999; $forecastData[52,1] = 5; $forecastData[52,2] = 1.022455007565251; $forecastData[52,3] = 10.94711643217653
$forecastData[53,0] = 45417.99999999999; $forecastData[53,1] = 6; $forecastData[53,2] = 0.2285215222913275; $forecastData[53,3] = 10.69942374476559
$forecastData[54,0] = 45424.99999999999; $forecastData[54,1] = 6; $forecastData[54,2] = 0.7086147202289022; $forecastData[54,3] = 10.67477736308166
$forecastData[55,0] = 45431.99999999999; $forecastData[55,1] = 6; $forecastData[55,2] = 0.567851704755615; $forecastData[55,3] = 10.90447883672405
$forecastData[56,0] = 45438.99999999999; $forecastData[56,1] = 6; $forecastData[56,2] = 0.7171206552538226; $forecastData[56,3] = 10.74181065714667
$wsForecast.Range("A2:D58").Value = $forecastData

# Reuse the date/time number format (the "ds" column) from an existing sheet
$wsWeekly.Range("A2").Copy()
$wsForecast.Range("A2:A58").PasteSpecial(-4122)
